$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''51.455.87'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '''3.107.10'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''389.01'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('D6').Value = '''104.10'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  -1.36%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '''37.19'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '''0.0859'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '''3.596.98'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '''18.59'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '''7.80'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').Value = '''3.116.88'
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').Value = '''10.64'
$ws.Range('E18').Value = '  +1.36%  '
$ws.Range('D19').Value = '''51.552.64'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '''3.27'
$ws.Range('E20').Value = '  +6.45%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = '''0.0₃0967'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').Value = '''70.33'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '''266.84'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').Value = '''3.19'
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').Value = '''8.08'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '''27.41'
$ws.Range('E27').Value = '  +4.16%  '
$ws.Range('D28').Value = '''7.20'
$ws.Range('E28').Value = '  -4.99%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '''0.165'
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').Value = '''10.44'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('D33').Value = '''35.88'
$ws.Range('E33').Value = '  +5.05%  '
$ws.Range('D34').Value = '''0.0476'
$ws.Range('E34').Value = '  +6.43%  '
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('D36').Value = '''50.07'
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('D39').Value = '''0.291'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('E40').Value = '  +1.47%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').Value = '''16.61'
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('D44').Value = '''3.83'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').Value = '''2.51'
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('D46').Value = '''22.33'
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('E47').Value = '  +4.24%  '
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').Value = '''2.080.47'
$ws.Range('E49').Value = '  +2.29%  '
$ws.Range('D50').Value = '''0.0332'
$ws.Range('E50').Value = '  +3.55%  '
$ws.Range('D51').Value = '''0.931'
$ws.Range('E51').Value = '  +19.06%  '
